# Generate Report for Handback
#
# A new handback-transform run for 6cbafaf2-1527-4d17-8609-c520de4665c1 has
# failed. The report generator re-emitted the three status-report rows that
# sit between "dbd0a76f..." and the ".localization-config" footer row:
# the freshly-processed item (6cbafaf2) now sorts first in that block (with
# its status flipped to "Handback transform failed" and a reset handback
# timestamp), pushing dc7a0273 and 50840c9f down by one row each, on every
# sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A6").Value = "6cbafaf2-1527-4d17-8609-c520de4665c1.md"
$ov.Range("B6").Value = "Handback transform failed"
$ov.Range("C6").Value = "Handback transform failed"

$ov.Range("A7").Value = "dc7a0273-dac9-469d-8c9c-361251acb6f8.md"
$ov.Range("B7").Value = "In Translation"
$ov.Range("C7").Value = "In Translation"

$ov.Range("A8").Value = "50840c9f-ebd6-46b3-ba37-fc06ee076493.md"
$ov.Range("B8").Value = "Ready for handoff"
$ov.Range("C8").Value = "Ready for handoff"

# ---- zh-cn sheet ---------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A6").Value = "6cbafaf2-1527-4d17-8609-c520de4665c1.md"
$zh.Range("B6").Value = "Handback transform failed"
$zh.Range("C6").Value = "6cbafaf2-1527-4d17-8609-c520de4665c1.29b749a075e8672c7bb2de352f8458798c8dc360.zh-cn.xlf"
$zh.Range("D6").Value = "2016-03-09 04:27:32"
$zh.Range("G6").Value = "0001-01-01 00:00:00"
$zh.Range("H6").Value = "Include"

$zh.Range("A7").Value = "dc7a0273-dac9-469d-8c9c-361251acb6f8.md"
$zh.Range("B7").Value = "In Translation"
$zh.Range("C7").Value = "dc7a0273-dac9-469d-8c9c-361251acb6f8.17a5d7b1dbb936cb3b2615b26a3f9d02f311d0f7.zh-cn.xlf"
$zh.Range("D7").Value = "2016-03-09 04:20:18"
$zh.Range("G7").Value = "0001-01-01 00:00:00"
$zh.Range("H7").Value = "Include"

$zh.Range("A8").Value = "50840c9f-ebd6-46b3-ba37-fc06ee076493.md"
$zh.Range("B8").Value = "Ready for handoff"
$zh.Range("C8").Value = "50840c9f-ebd6-46b3-ba37-fc06ee076493.3aafc46c3c43d4a6668076903881bc9086c03c65.zh-cn.xlf"
$zh.Range("D8").Value = "2016-03-09 04:27:32"
$zh.Range("G8").Value = "0001-01-01 00:00:00"
$zh.Range("H8").Value = "Include"

# ---- de-de sheet ---------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A6").Value = "6cbafaf2-1527-4d17-8609-c520de4665c1.md"
$de.Range("B6").Value = "Handback transform failed"
$de.Range("C6").Value = "6cbafaf2-1527-4d17-8609-c520de4665c1.29b749a075e8672c7bb2de352f8458798c8dc360.de-de.xlf"
$de.Range("D6").Value = "2016-03-09 04:27:35"
$de.Range("G6").Value = "0001-01-01 00:00:00"
$de.Range("H6").Value = "Include"

$de.Range("A7").Value = "dc7a0273-dac9-469d-8c9c-361251acb6f8.md"
$de.Range("B7").Value = "In Translation"
$de.Range("C7").Value = "dc7a0273-dac9-469d-8c9c-361251acb6f8.17a5d7b1dbb936cb3b2615b26a3f9d02f311d0f7.de-de.xlf"
$de.Range("D7").Value = "2016-03-09 04:20:21"
$de.Range("G7").Value = "0001-01-01 00:00:00"
$de.Range("H7").Value = "Include"

$de.Range("A8").Value = "50840c9f-ebd6-46b3-ba37-fc06ee076493.md"
$de.Range("B8").Value = "Ready for handoff"
$de.Range("C8").Value = "50840c9f-ebd6-46b3-ba37-fc06ee076493.3aafc46c3c43d4a6668076903881bc9086c03c65.de-de.xlf"
$de.Range("D8").Value = "2016-03-09 04:27:35"
$de.Range("G8").Value = "0001-01-01 00:00:00"
$de.Range("H8").Value = "Include"
